# Revert commit c1bab2bfb1886b32f82a259f7ace6b03368253a6
# ("Descrição dos casos de uso - Faculdade")
#
# 1. Bump every table's left indentation (w:tblInd) and every left cell
#    margin (w:tblCellMar / w:tcMar) by 5 dxa (= 0.25 pt) across all 7
#    use-case tables.
# 2. Drop 3 of the 4 trailing empty "Corpodetexto" (Body Text) paragraphs
#    after the last table, keeping just one.
# 3. Delete the unused ListLabel10..ListLabel18 character styles.
# 4. Reset the "Normal" style's font color back to Automatic.

$d = $word.ActiveDocument

# --- 1. Table indentation / cell margins -----------------------------
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $t = $d.Tables.Item($ti)

    $t.Rows.LeftIndent = $t.Rows.LeftIndent + 0.25
    $t.LeftPadding = $t.LeftPadding + 0.25

    foreach ($r in $t.Rows) {
        foreach ($c in $r.Cells) {
            $c.LeftPadding = $c.LeftPadding + 0.25
        }
    }
}

# --- 2. Trim trailing empty paragraphs after the last table -----------
$n = $d.Paragraphs.Count
$d.Paragraphs.Item($n - 1).Range.Delete()
$d.Paragraphs.Item($n - 2).Range.Delete()
$d.Paragraphs.Item($n - 3).Range.Delete()

# --- 3. Remove unused ListLabel10..ListLabel18 character styles -------
for ($i = 10; $i -le 18; $i++) {
    $s = $d.Styles.Item("ListLabel" + $i)
    $s.Delete()
}

# --- 4. Normal style font color back to Automatic ---------------------
$normal = $d.Styles.Item("Normal")
$normal.Font.Color = -16777216
